$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1521.238
$ws.Range("J32").Value = 1541.4445
$ws.Range("L32").Value = 1541.4445
$ws.Range("N32").Value = -2193.4445
$ws.Range("H64").Value = 3492.375
$ws.Range("I64").Value = 3496.9092
$ws.Range("J64").Value = 3488.5386
$ws.Range("K64").Value = 3496.9092
$ws.Range("L64").Value = 3488.5386
$ws.Range("M64").Value = -3248.9092
$ws.Range("N64").Value = -3984.5386
$ws.Range("H67").Value = 3492.375
$ws.Range("I67").Value = 3496.9092
$ws.Range("J67").Value = 3488.5386
$ws.Range("K67").Value = 3496.9092
$ws.Range("L67").Value = 3488.5386
$ws.Range("M67").Value = -2638.9092
$ws.Range("N67").Value = -5204.5386
$ws.Range("H74").Value = 10714.286
$ws.Range("I74").Value = 12750
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 12750
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -11814
$ws.Range("N74").Value = -9872
$ws.Range("H77").Value = 10714.286
$ws.Range("I77").Value = 12750
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 63750
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -59070
$ws.Range("N77").Value = -49360
$ws.Range("H125").Value = 1502
$ws.Range("I125").Value = 1760.5
$ws.Range("K125").Value = 15844.5
$ws.Range("M125").Value = -13384.5
$ws.Range("H129").Value = 846.3461
$ws.Range("J129").Value = 885.4789
$ws.Range("L129").Value = 2656.4367
$ws.Range("N129").Value = -12656.4367
$ws.Range("H132").Value = 10757161
$ws.Range("I132").Value = 11496475
$ws.Range("J132").Value = 37100
$ws.Range("K132").Value = 34489425
$ws.Range("L132").Value = 111300
$ws.Range("M132").Value = -34486895
$ws.Range("N132").Value = -116360
$ws.Range("H135").Value = 349
$ws.Range("I135").Value = 273.83334
$ws.Range("K135").Value = 2464.50006
$ws.Range("M135").Value = 70.4999399999997
$ws.Range("H137").Value = 1121.862
$ws.Range("I137").Value = 1017.75
$ws.Range("K137").Value = 3053.25
$ws.Range("M137").Value = -503.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4150.5435
$ws.Range("I32").Value = 3934.1538
$ws.Range("J32").Value = 5356.143
$ws.Range("K32").Value = 3934.1538
$ws.Range("L32").Value = 5356.143
$ws.Range("M32").Value = -3647.1538
$ws.Range("N32").Value = -5930.143
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H61").Value = 1881.3572
$ws.Range("I61").Value = 1394.4546
$ws.Range("J61").Value = 3666.6667
$ws.Range("K61").Value = 1394.4546
$ws.Range("L61").Value = 3666.6667
$ws.Range("M61").Value = -1182.4546
$ws.Range("N61").Value = -4090.6667
$ws.Range("H74").Value = 988.1539
$ws.Range("J74").Value = 1569
$ws.Range("L74").Value = 1569
$ws.Range("N74").Value = -3317
$ws.Range("H77").Value = 988.1539
$ws.Range("J77").Value = 1569
$ws.Range("L77").Value = 7845
$ws.Range("N77").Value = -16581
$ws.Range("H101").Value = 35602
$ws.Range("J101").Value = 35602
$ws.Range("L101").Value = 35602
$ws.Range("N101").Value = -42092
$ws.Range("H119").Value = 22999.666
$ws.Range("J119").Value = 22999.666
$ws.Range("L119").Value = 22999.666
$ws.Range("N119").Value = -32675.666
$ws.Range("H125").Value = 44800
$ws.Range("J125").Value = 44800
$ws.Range("L125").Value = 44800
$ws.Range("N125").Value = -54640
$ws.Range("H132").Value = 2425.2632
$ws.Range("I132").Value = 2073
$ws.Range("J132").Value = 4304
$ws.Range("K132").Value = 6219
$ws.Range("L132").Value = 12912
$ws.Range("M132").Value = -3689
$ws.Range("N132").Value = -17972
$ws.Range("H136").Value = 1881.3572
$ws.Range("I136").Value = 1394.4546
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 4183.3638
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -1633.3638
$ws.Range("N136").Value = -16100.0001
$ws.Range("H139").Value = 30659.75
$ws.Range("J139").Value = 30659.75
$ws.Range("L139").Value = 30659.75
$ws.Range("N139").Value = -40939.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12175.333
$ws.Range("I134").Value = 1872.1428
$ws.Range("J134").Value = 26599.8
$ws.Range("K134").Value = 5616.428400000001
$ws.Range("L134").Value = 79799.39999999999
$ws.Range("M134").Value = -3081.428400000001
$ws.Range("N134").Value = -84869.39999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1229.2727
$ws.Range("I58").Value = 1054.9474
$ws.Range("J58").Value = 2333.3333
$ws.Range("K58").Value = 1054.9474
$ws.Range("L58").Value = 2333.3333
$ws.Range("M58").Value = -851.9474
$ws.Range("N58").Value = -2739.3333
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 2203.7896
$ws.Range("I132").Value = 1598.1428
$ws.Range("J132").Value = 3899.6
$ws.Range("K132").Value = 4794.428400000001
$ws.Range("L132").Value = 11698.8
$ws.Range("M132").Value = -2264.428400000001
$ws.Range("N132").Value = -16758.8
$ws.Range("H134").Value = 19609166
$ws.Range("I134").Value = 30304186
$ws.Range("J134").Value = 1627.8334
$ws.Range("K134").Value = 90912558
$ws.Range("L134").Value = 4883.5002
$ws.Range("M134").Value = -90910023
$ws.Range("N134").Value = -9953.5002
$ws.Range("H136").Value = 1229.2727
$ws.Range("I136").Value = 1054.9474
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 3164.8422
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -614.8422
$ws.Range("N136").Value = -12099.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2012.5714
$ws.Range("J69").Value = 2166.5264
$ws.Range("L69").Value = 6499.5792
$ws.Range("N69").Value = -8121.5792
$ws.Range("H72").Value = 2012.5714
$ws.Range("J72").Value = 2166.5264
$ws.Range("L72").Value = 19498.7376
$ws.Range("N72").Value = -27610.7376
$ws.Range("H128").Value = 50000
$ws.Range("I128").Value = 50000
$ws.Range("K128").Value = 150000
$ws.Range("M128").Value = -145020
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 2584.5
$ws.Range("I132").Value = 2001.7778
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 6005.3334
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -3475.3334
$ws.Range("N132").Value = -18057.9995
$ws.Range("H133").Value = 42134.625
$ws.Range("J133").Value = 42134.625
$ws.Range("L133").Value = 42134.625
$ws.Range("N133").Value = -52254.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1551.3043
$ws.Range("I22").Value = 1540.0588
$ws.Range("J22").Value = 1583.1666
$ws.Range("K22").Value = 1540.0588
$ws.Range("L22").Value = 1583.1666
$ws.Range("M22").Value = -1245.0588
$ws.Range("N22").Value = -2173.1666
$ws.Range("H27").Value = 1551.3043
$ws.Range("I27").Value = 1540.0588
$ws.Range("J27").Value = 1583.1666
$ws.Range("K27").Value = 1540.0588
$ws.Range("L27").Value = 1583.1666
$ws.Range("M27").Value = -1433.0588
$ws.Range("N27").Value = -1797.1666
$ws.Range("H40").Value = 3041
$ws.Range("I40").Value = 2966.6667
$ws.Range("J40").Value = 3152.5
$ws.Range("K40").Value = 2966.6667
$ws.Range("L40").Value = 3152.5
$ws.Range("M40").Value = -2830.6667
$ws.Range("N40").Value = -3424.5
$ws.Range("H68").Value = 1824.3077
$ws.Range("I68").Value = 1551.8
$ws.Range("J68").Value = 2732.6667
$ws.Range("K68").Value = 1551.8
$ws.Range("L68").Value = 2732.6667
$ws.Range("M68").Value = -802.8
$ws.Range("N68").Value = -4230.6667
$ws.Range("H71").Value = 1824.3077
$ws.Range("I71").Value = 1551.8
$ws.Range("J71").Value = 2732.6667
$ws.Range("K71").Value = 7759
$ws.Range("L71").Value = 13663.3335
$ws.Range("M71").Value = -4015
$ws.Range("N71").Value = -21151.3335
$ws.Range("H122").Value = 27780820
$ws.Range("I122").Value = 41669732
$ws.Range("K122").Value = 125009196
$ws.Range("M122").Value = -125006746
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H127").Value = 35357.5
$ws.Range("J127").Value = 35357.5
$ws.Range("L127").Value = 35357.5
$ws.Range("N127").Value = -45277.5
$ws.Range("H136").Value = 1966.6364
$ws.Range("I136").Value = 1766.75
$ws.Range("K136").Value = 5300.25
$ws.Range("M136").Value = -2750.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 125008750
$ws.Range("I62").Value = 166676670
$ws.Range("J62").Value = 4999
$ws.Range("K62").Value = 166676670
$ws.Range("L62").Value = 4999
$ws.Range("M62").Value = -166676046
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 125008750
$ws.Range("I65").Value = 166676670
$ws.Range("J65").Value = 4999
$ws.Range("K65").Value = 833383350
$ws.Range("L65").Value = 24995
$ws.Range("M65").Value = -833380230
$ws.Range("N65").Value = -31235
$ws.Range("H122").Value = 23640042
$ws.Range("I122").Value = 23640042
$ws.Range("K122").Value = 70920126
$ws.Range("M122").Value = -70917676
$ws.Range("H132").Value = 5026.7334
$ws.Range("I132").Value = 4619.4546
$ws.Range("K132").Value = 13858.3638
$ws.Range("M132").Value = -11328.3638
$ws.Range("H136").Value = 1003
$ws.Range("I136").Value = 804.1
$ws.Range("J136").Value = 1666
$ws.Range("K136").Value = 2412.3
$ws.Range("L136").Value = 4998
$ws.Range("M136").Value = 137.6999999999998
$ws.Range("N136").Value = -10098
